$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 with the next order entry (continuing the pattern from rows 2-3).
# Reuse the same date number format as the cells above (A2:A3) so the new
# cell picks up the existing style instead of minting a new one.
$dateFormat = $ws.Range("A3").NumberFormat
$ws.Range("A4").Value = 44486
$ws.Range("A4").NumberFormat = $dateFormat
$ws.Range("B4").Value = "URC-2021-12-19"
$ws.Range("C4").Value = "9W"

# Move the active selection to C10, matching the post-edit view state
$ws.Range("C10").Select()
